# TALLER 2 PRESENTACIÓN - añadida comparacion con Corgi
#
# A new slide ("¿Comparación con Corgi?") is inserted at position 4
# (right after "¿QUÉ ES UN EVENT EN C#?" and before "¿QUÉ ES UN SINGLETON?"),
# using the "Título y objetos" (Title and Content) layout - the same layout
# used by its neighboring slides.

$p = $ppt.ActivePresentation

# "Título y objetos" is CustomLayout index 2 on the deck's single slide master.
$layout = $p.SlideMaster.CustomLayouts.Item(2)

# Insert the new slide at slide position 4.
$newSlide = $p.Slides.AddSlide(4, $layout)

# --- Title placeholder -------------------------------------------------
$titleShape = $newSlide.Shapes.Item(1)
$titleShape.Name = "Título 1"

$titleRange = $titleShape.TextFrame.TextRange
$titleRange.Text = "Comparacion"
$titleRange.InsertAfter(" con ") | Out-Null
$titleRange.InsertAfter("Corgi") | Out-Null

# --- Body / content placeholder -----------------------------------------
$bodyShape = $newSlide.Shapes.Item(2)
$bodyShape.Name = "Marcador de contenido 2"

$bodyRange = $bodyShape.TextFrame.TextRange
$bodyRange.Text = "El concepto de eventos en el "
$bodyRange.InsertAfter("Corgi") | Out-Null
$bodyRange.InsertAfter(" es el mismo que en c#, un suceso al cual se le pueden suscribir sujetos y se les notifica si el evento ocurre. También usa delegados como los eventos clásicos del c#, la diferencia es que ") | Out-Null
$bodyRange.InsertAfter("Corgi") | Out-Null
$bodyRange.InsertAfter(" tiene simplificado el proceso y con varios eventos predefinidos de ") | Out-Null
$bodyRange.InsertAfter("sucesos comunes ") | Out-Null
$bodyRange.InsertAfter("en ") | Out-Null
$bodyRange.InsertAfter("juegos 2D.") | Out-Null
